$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update data cell A2: "Principal" -> "PI"
$ws.Range("A2").Value = "PI"

# Update header cell A1: "Role (Principal/Sub)" -> "Role (PI/Sub I)"
$ws.Range("A1").Value = "Role (PI/Sub I)"

# Update data cell A3: "Principal" -> "Sub I"
$ws.Range("A3").Value = "Sub I"
